# SW requirements were updated.
# Rebuild the "wishlist" requirements table: add Type/Status columns and
# expand the requirement rows (incl. new sub-requirement / definition rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SWRS")
$lo = $ws.ListObjects.Item("wishlist")

# ---------------------------------------------------------------------------
# 1. Grow the table to B2:E9 (2 new columns, 5 new rows) before touching data
#    so every cell we write below already belongs to the ListObject.
# ---------------------------------------------------------------------------
$lo.Resize($ws.Range("B2:E9"))

# Table used to show autofilter drop-down buttons on its two columns; the
# rebuilt table no longer shows them.
$lo.ShowAutoFilter = $false

# ---------------------------------------------------------------------------
# 2. Formatting first (copy from the existing, already-styled cells), then
#    fill in text - this keeps every cell on the style the real workbook
#    uses instead of the generic "Normal" style new range cells start on.
# ---------------------------------------------------------------------------

# Header cells D2:E2 look like the existing header cells B2:C2.
$ws.Range("C2").Copy() | Out-Null
$ws.Range("D2:E2").PasteSpecial(-4122) | Out-Null

# Row 3's Type/Status cells (D3:E3) reuse the existing data-row look
# (font, left/center alignment) but centered - this is the one new cell
# style the workbook gains.
$ws.Range("C3").Copy() | Out-Null
$ws.Range("D3:E3").PasteSpecial(-4122) | Out-Null
$ws.Range("D3:E3").HorizontalAlignment = -4108

# Row 3 (now fully styled, B3:E3) becomes the template for every other
# data row (4-9), both new and pre-existing.
$ws.Range("B3:E3").Copy() | Out-Null
$ws.Range("B4:E9").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false

# Row heights for the newly added rows match the other data rows.
$ws.Range("B4:B9").RowHeight = 18.95

# ---------------------------------------------------------------------------
# 3. Header row text.
# ---------------------------------------------------------------------------
$ws.Range("B2").Value = "Requirement ID"
$ws.Range("C2").Value = "Requirement Specification"
$ws.Range("D2").Value = "Type"
$ws.Range("E2").Value = "Status"

# ---------------------------------------------------------------------------
# 4. Data rows.
# ---------------------------------------------------------------------------
$rows = @(
    @{ Id = "SWRS_BRMTR_001";   Spec = "The result files(s) shall be archived after finishing the test.";                            Type = "Requirement"; Status = "Open" }
    @{ Id = "SWRS_BRMTR_001_1"; Spec = "Archive name must be unique.";                                                                 Type = "Definition";  Status = "Open" }
    @{ Id = "SWRS_BRMTR_002";   Spec = "Archive containing test results shall be saved on local hard drive.";                         Type = "Requirement"; Status = "Open" }
    @{ Id = "SWRS_BRMTR_003";   Spec = "The user shall be promted to e-mail the test archive to borbalabc@gmail.com";                 Type = "Requirement"; Status = "Open" }
    @{ Id = "SWRS_BRMTR_004";   Spec = "User information and test result shall be stored in a text file.";                            Type = "Requirement"; Status = "Open" }
    @{ Id = "SWRS_BRMTR_004_1"; Spec = "The name of the text file must be a unique ID.";                                               Type = "Definition";  Status = "Open" }
    @{ Id = "SWRS_BRMTR_005";   Spec = "A unique ID shall be calculated for each user, based on their personal information.";          Type = "Requirement"; Status = "Open" }
)

$r = 3
foreach ($row in $rows) {
    $ws.Cells.Item($r, 2).Value = $row.Id
    $ws.Cells.Item($r, 3).Value = $row.Spec
    $ws.Cells.Item($r, 4).Value = $row.Type
    $ws.Cells.Item($r, 5).Value = $row.Status
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# 5. Column headers on the ListObject follow the header-row text written
#    above; make the names explicit too.
# ---------------------------------------------------------------------------
$lo.ListColumns.Item(1).Name = "Requirement ID"
$lo.ListColumns.Item(2).Name = "Requirement Specification"
$lo.ListColumns.Item(3).Name = "Type"
$lo.ListColumns.Item(4).Name = "Status"

# ---------------------------------------------------------------------------
# 6. Column widths: ID column widened (now holds longer "_1" ids) and the
#    new Type column sized to fit "Requirement"/"Definition".
# ---------------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 18.57
$ws.Columns.Item(4).ColumnWidth = 12.43

# ---------------------------------------------------------------------------
# 7. Selection left where the author ended up after the edit.
# ---------------------------------------------------------------------------
$ws.Range("E7").Select() | Out-Null
